{"js": "// Update the worksheet date header and all of the two-digit multiplication\n// prompts in the table. Each old prompt string is unique in the document, so\n// a straightforward search-and-replace keyed off the old text is safe.\nconst replacements = [\n  [\"2024-01-21 Sunday\", \"2024-01-22 Monday\"],\n  [\"96\u00d757=\", \"94\u00d772=\"],\n  [\"38\u00d723=\", \"51\u00d749=\"],\n  [\"26\u00d740=\", \"49\u00d742=\"],\n  [\"28\u00d732=\", \"96\u00d760=\"],\n  [\"88\u00d715=\", \"78\u00d752=\"],\n  [\"92\u00d752=\", \"11\u00d734=\"],\n  [\"66\u00d724=\", \"25\u00d783=\"],\n  [\"85\u00d791=\", \"17\u00d766=\"],\n  [\"92\u00d779=\", \"19\u00d770=\"],\n  [\"68\u00d729=\", \"63\u00d765=\"],\n  [\"65\u00d758=\", \"69\u00d743=\"],\n  [\"46\u00d733=\", \"47\u00d797=\"],\n  [\"66\u00d725=\", \"32\u00d717=\"],\n  [\"74\u00d767=\", \"32\u00d746=\"],\n  [\"64\u00d752=\", \"75\u00d781=\"],\n  [\"15\u00d734=\", \"46\u00d714=\"],\n  [\"28\u00d792=\", \"90\u00d755=\"],\n  [\"68\u00d754=\", \"26\u00d750=\"],\n  [\"45\u00d739=\", \"47\u00d744=\"],\n  [\"89\u00d771=\", \"90\u00d735=\"],\n  [\"30\u00d745=\", \"66\u00d748=\"],\n  [\"27\u00d793=\", \"19\u00d755=\"],\n  [\"45\u00d719=\", \"24\u00d773=\"],\n  [\"34\u00d732=\", \"57\u00d758=\"],\n  [\"41\u00d731=\", \"19\u00d775=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date header and all of the two-digit multiplication\n# prompts in the table. Each old prompt string is unique in the document, so\n# a Find/Replace (ReplaceAll) pass keyed off the old text is safe for each.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-01-21 Sunday\", \"2024-01-22 Monday\"),\n  @(\"96\u00d757=\", \"94\u00d772=\"),\n  @(\"38\u00d723=\", \"51\u00d749=\"),\n  @(\"26\u00d740=\", \"49\u00d742=\"),\n  @(\"28\u00d732=\", \"96\u00d760=\"),\n  @(\"88\u00d715=\", \"78\u00d752=\"),\n  @(\"92\u00d752=\", \"11\u00d734=\"),\n  @(\"66\u00d724=\", \"25\u00d783=\"),\n  @(\"85\u00d791=\", \"17\u00d766=\"),\n  @(\"92\u00d779=\", \"19\u00d770=\"),\n  @(\"68\u00d729=\", \"63\u00d765=\"),\n  @(\"65\u00d758=\", \"69\u00d743=\"),\n  @(\"46\u00d733=\", \"47\u00d797=\"),\n  @(\"66\u00d725=\", \"32\u00d717=\"),\n  @(\"74\u00d767=\", \"32\u00d746=\"),\n  @(\"64\u00d752=\", \"75\u00d781=\"),\n  @(\"15\u00d734=\", \"46\u00d714=\"),\n  @(\"28\u00d792=\", \"90\u00d755=\"),\n  @(\"68\u00d754=\", \"26\u00d750=\"),\n  @(\"45\u00d739=\", \"47\u00d744=\"),\n  @(\"89\u00d771=\", \"90\u00d735=\"),\n  @(\"30\u00d745=\", \"66\u00d748=\"),\n  @(\"27\u00d793=\", \"19\u00d755=\"),\n  @(\"45\u00d719=\", \"24\u00d773=\"),\n  @(\"34\u00d732=\", \"57\u00d758=\"),\n  @(\"41\u00d731=\", \"19\u00d775=\")\n)\n\nforeach ($p in $pairs) {\n  $rng = $d.Content\n  $rng.Find.Execute($p[0], $false, $false, $false, $false, $false, $true, 1, $false, $p[1], 2)\n}\n"}
